# Actualización automática 2025-06-18 13:35:10
# Update June ("junio") sales figure for LOZANO MOLINA TITO / RENOVA&DISEÑA S.A.
# on the "240X80 PORCELANATO" product group, and propagate the change through
# the dependent totals on the other sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO": column D (240X80 PORCELANATO), row 19 ---
$wsVentasPorGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentasPorGrupo.Range("D19").Value = 366.34

# --- Sheet "VENTA MENSUAL": column F (junio), row 19, and the total row 29 ---
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsVentaMensual.Range("F19").Value = 366.34
$wsVentaMensual.Range("F29").Value = 6301.41

# --- Sheet "CUMPLIMIENTO MENSUAL": row 3 (240X80 PORCELANATO) and row 19 (TOTAL) ---
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumplimiento.Range("D3").Value = 2528.06
$wsCumplimiento.Range("E3").Value = 592.0545000000002
$wsCumplimiento.Range("F3").Value = 0.8102459060396662

$wsCumplimiento.Range("D19").Value = 12689.86
$wsCumplimiento.Range("E19").Value = 10810.14093005039
$wsCumplimiento.Range("F19").Value = 0.5399940211820574
